$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 39.18403521750415
$ws.Range("C2").Value = 38.48091759558707
$ws.Range("D2").Value = 41.13242413897014
$ws.Range("E2").Value = 42.0807718115307

$ws.Range("B3").Value = 50.4052041242598
$ws.Range("C3").Value = 50.29720458794134
$ws.Range("D3").Value = 50.22638394165971
$ws.Range("E3").Value = 50.27275885130553

$ws.Range("B4").Value = 98.77004542256238
$ws.Range("C4").Value = 98.70795095373275
$ws.Range("D4").Value = 98.8258534131719
$ws.Range("E4").Value = 98.89918668649058

$ws.Range("B5").Value = 98.95855301637833
$ws.Range("C5").Value = 98.95681476683551
$ws.Range("D5").Value = 98.88779824785216
$ws.Range("E5").Value = 98.90346527239799

$ws.Range("B6").Value = 98.54094725838611
$ws.Range("C6").Value = 98.48428624236892
$ws.Range("D6").Value = 98.49223775236969
$ws.Range("E6").Value = 98.44336705398108

$ws.Range("B7").Value = 98.02445065369528
$ws.Range("C7").Value = 98.02309691807979
$ws.Range("D7").Value = 98.0631857166878
$ws.Range("E7").Value = 98.00601175557821

$ws.Range("B8").Value = 97.59316819264832
$ws.Range("C8").Value = 97.49565773388127
$ws.Range("D8").Value = 97.62472004143031
$ws.Range("E8").Value = 97.54299680541318

$ws.Range("B9").Value = 96.26427673980157
$ws.Range("C9").Value = 96.28451487404142
$ws.Range("D9").Value = 96.28558978656022
$ws.Range("E9").Value = 96.28849438021088
